# Add a new user row ("Anna Yukimi Yamada" / "iuquimi" / "123456" / "BMW")
# to the query table on the "query" sheet, growing the table/range from
# A1:D27 to A1:D28.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("query")

$newRow = 28

# Match the formatting (text number format, same as the rest of the table
# body) before writing values, so e.g. the numeric-looking password
# "123456" is stored as text rather than being reinterpreted as a number.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 2).NumberFormat = "@"
$ws.Cells.Item($newRow, 3).NumberFormat = "@"
$ws.Cells.Item($newRow, 4).NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "Anna Yukimi Yamada"
$ws.Cells.Item($newRow, 2).Value = "iuquimi"
$ws.Cells.Item($newRow, 3).Value = "123456"
$ws.Cells.Item($newRow, 4).Value = "BMW"

# Resize the query table (and its autofilter) to include the new row.
$table = $ws.ListObjects.Item("Table_query")
$table.Resize($ws.Range("A1:D$newRow"))

# Update the hidden workbook-level defined name that mirrors the table range.
$wb.Names.Item("query").RefersTo = "=query!`$A`$1:`$D`$$newRow"

# Move the active selection to the next empty row, as Excel does after
# typing a new row of data into the table.
$ws.Range("D" + ($newRow + 1)).Select()
